$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new line rows ("line7", "line8") are inserted right after "line6"
# (pushing the former extr1..extr8 rows down from 8-15 to 10-17), and the
# C/D/E (from_bus/to_bus/in_service) values for rows 8-17 are updated to
# match the new contingency results.

# First, copy the formatting of the last existing data row (A15, which
# carries style index 1) down onto the two brand-new rows 16 and 17 so
# that column A keeps the same bold/bordered/centered style.
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16:A17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$data = @(
    @(8,  6, "line7", 14, 11, $true),
    @(9,  7, "line8", 16, 9,  $true),
    @(10, 8, "extr1", 5,  12, $true),
    @(11, 9, "extr2", 5,  9,  $true),
    @(12, 10, "extr3", 10, 11, $false),
    @(13, 11, "extr4", 7,  8,  $true),
    @(14, 12, "extr5", 9,  11, $false),
    @(15, 13, "extr6", 7,  11, $false),
    @(16, 14, "extr7", 5,  7,  $true),
    @(17, 15, "extr8", 8,  5,  $false)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
